# Update cryptocurrency price/volume data per latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.759.84"
$ws.Range("E2").Value = '  +3.39%  '

# Row 3
$ws.Range("D3").Value = "'2.446.46"
$ws.Range("E3").Value = '  +2.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").Value = "'575.63"
$ws.Range("E5").Value = '  +2.56%  '

# Row 6
$ws.Range("D6").Value = "'145.76"
$ws.Range("E6").Value = '  +3.19%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
$ws.Range("D9").Value = "'2.444.41"
$ws.Range("E9").Value = '  +1.71%  '

# Row 10
$ws.Range("E10").Value = '  +3.02%  '

# Row 11
$ws.Range("E11").Value = '  +2.55%  '

# Row 12
$ws.Range("E12").Value = '  +1.32%  '

# Row 13
$ws.Range("E13").Value = '  +2.40%  '

# Row 14
$ws.Range("D14").Value = "'28.20"
$ws.Range("E14").Value = '  +7.27%  '

# Row 15
$ws.Range("E15").Value = '  +6.13%  '

# Row 16
$ws.Range("D16").Value = "'2.889.60"
$ws.Range("E16").Value = '  +3.13%  '

# Row 17
$ws.Range("D17").Value = "'62.559.83"
$ws.Range("E17").Value = '  +3.76%  '

# Row 18
$ws.Range("D18").Value = "'2.438.19"
$ws.Range("E18").Value = '  +1.42%  '

# Row 19
$ws.Range("D19").Value = "'7.91"
$ws.Range("E19").Value = '  -2.26%  '

# Row 20
$ws.Range("E20").Value = '  +2.97%  '

# Row 21
$ws.Range("D21").Value = "'329.56"
$ws.Range("E21").Value = '  +1.62%  '

# Row 22
$ws.Range("E22").Value = '  +1.14%  '

# Row 23
$ws.Range("D23").Value = "'2.04"
$ws.Range("E23").Value = '  +9.45%  '

# Row 24
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = "'66.05"
$ws.Range("E25").Value = '  +1.95%  '

# Row 26
$ws.Range("B26").Value = 'BabyDogeCoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D26").Value = "'0.0₆0618"
$ws.Range("E26").Value = '  +122.13%  '

# Row 27
$ws.Range("D27").Value = "'651.92"
$ws.Range("E27").Value = '  +12.80%  '

# Row 28
$ws.Range("E28").Value = '  +17.61%  '

# Row 29
$ws.Range("E29").Value = '  +4.93%  '

# Row 30
$ws.Range("D30").Value = "'0.0₃0990"
$ws.Range("E30").Value = '  +5.86%  '

# Row 31
$ws.Range("D31").Value = "'2.567.08"
$ws.Range("E31").Value = '  +2.13%  '

# Row 32
$ws.Range("E32").Value = '  +8.97%  '

# Row 33
$ws.Range("D33").Value = "'8.22"
$ws.Range("E33").Value = '  +2.25%  '

# Row 34
$ws.Range("E34").Value = '  +3.35%  '

# Row 36
$ws.Range("E36").Value = '  +2.41%  '

# Row 38
$ws.Range("E38").Value = '  +3.67%  '

# Row 39
$ws.Range("D39").Value = "'5.51"
$ws.Range("E39").Value = '  +6.85%  '

# Row 40
$ws.Range("D40").Value = "'0.375"
$ws.Range("E40").Value = '  +1.01%  '

# Row 41
$ws.Range("D41").Value = "'152.58"
$ws.Range("E41").Value = '  +0.15%  '

# Row 42
$ws.Range("E42").Value = '  +2.48%  '

# Row 43
$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = '  +9.00%  '

# Row 44
$ws.Range("E44").Value = '  +5.22%  '

# Row 45
$ws.Range("D45").Value = "'42.39"
$ws.Range("E45").Value = '  +1.81%  '

# Row 46
$ws.Range("E46").Value = '  +0.00%  '

# Row 47
$ws.Range("E47").Value = '  +27.43%  '

# Row 48
$ws.Range("D48").Value = "'145.50"
$ws.Range("E48").Value = '  +3.00%  '

# Row 49
$ws.Range("E49").Value = '  +3.44%  '

# Row 50
$ws.Range("D50").Value = "'20.77"
$ws.Range("E50").Value = '  +7.17%  '

# Row 51
$ws.Range("E51").Value = '  +2.83%  '
